# chore: update Sheets via scheduled runner
#
# Refreshes the market-board snapshot columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> columns
# H through N) on every per-job worksheet (ALC, ARM, BSM, CRP, CUL, GSM,
# LTW, WVR) with newly scraped values.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 996.6667
$ws.Range("I9").Value = 996.6667
$ws.Range("K9").Value = 996.6667
$ws.Range("M9").Value = -827.6667
$ws.Range("H40").Value = 3499
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3499
$ws.Range("K40").Value = 0
$ws.Range("N40").Value = -3849
$ws.Range("H132").Value = 2602.2856
$ws.Range("I132").Value = 1951.0435
$ws.Range("J132").Value = 5598
$ws.Range("K132").Value = 5853.1305
$ws.Range("L132").Value = 16794
$ws.Range("M132").Value = -3323.1305
$ws.Range("N132").Value = -21854
$ws.Range("H137").Value = 2534.484
$ws.Range("I137").Value = 2291.9333
$ws.Range("K137").Value = 6875.7999
$ws.Range("M137").Value = -4325.7999
$ws.Range("L40").Value = 3499
$ws.Range("M40").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3177.6753
$ws.Range("I32").Value = 3028.8933
$ws.Range("K32").Value = 3028.8933
$ws.Range("M32").Value = -2741.8933
$ws.Range("H57").Value = 5499.6665
$ws.Range("I57").Value = 5499.6665
$ws.Range("K57").Value = 5499.6665
$ws.Range("M57").Value = -5015.6665
$ws.Range("H61").Value = 5519.4385
$ws.Range("I61").Value = 2733.182
$ws.Range("K61").Value = 2733.182
$ws.Range("M61").Value = -2521.182
$ws.Range("H63").Value = 4999.5
$ws.Range("I63").Value = 4999.5
$ws.Range("K63").Value = 4999.5
$ws.Range("M63").Value = -4313.5
$ws.Range("H66").Value = 4999.5
$ws.Range("I66").Value = 4999.5
$ws.Range("K66").Value = 24997.5
$ws.Range("M66").Value = -21565.5
$ws.Range("H102").Value = 1303
$ws.Range("I102").Value = 1377.1818
$ws.Range("J102").Value = 1139.8
$ws.Range("K102").Value = 1377.1818
$ws.Range("L102").Value = 1139.8
$ws.Range("M102").Value = 244.8181999999999
$ws.Range("N102").Value = -4383.8
$ws.Range("H122").Value = 20712.076
$ws.Range("I122").Value = 27669.625
$ws.Range("J122").Value = 9580
$ws.Range("K122").Value = 83008.875
$ws.Range("L122").Value = 28740
$ws.Range("M122").Value = -80558.875
$ws.Range("N122").Value = -33640
$ws.Range("H126").Value = 5245.3335
$ws.Range("I126").Value = 5245.3335
$ws.Range("K126").Value = 15736.0005
$ws.Range("M126").Value = -13266.0005
$ws.Range("H136").Value = 5519.4385
$ws.Range("I136").Value = 2733.182
$ws.Range("K136").Value = 8199.545999999998
$ws.Range("M136").Value = -5649.545999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6174332
$ws.Range("I20").Value = 8334949.5
$ws.Range("K20").Value = 8334949.5
$ws.Range("M20").Value = -8334702.5
$ws.Range("H86").Value = 35717092
$ws.Range("I86").Value = 14708316
$ws.Range("J86").Value = 55558716
$ws.Range("K86").Value = 14708316
$ws.Range("L86").Value = 55558716
$ws.Range("M86").Value = -14707193
$ws.Range("N86").Value = -55560962
$ws.Range("H89").Value = 35717092
$ws.Range("I89").Value = 14708316
$ws.Range("J89").Value = 55558716
$ws.Range("K89").Value = 73541580
$ws.Range("L89").Value = 277793580
$ws.Range("M89").Value = -73535964
$ws.Range("N89").Value = -277804812
$ws.Range("H113").Value = 5036
$ws.Range("I113").Value = 5036
$ws.Range("K113").Value = 5036
$ws.Range("M113").Value = -2866
$ws.Range("H134").Value = 5499.067
$ws.Range("I134").Value = 2044.8276
$ws.Range("K134").Value = 6134.4828
$ws.Range("M134").Value = -3599.4828

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8988.977000000001
$ws.Range("I31").Value = 4072.7222
$ws.Range("J31").Value = 12528.68
$ws.Range("K31").Value = 4072.7222
$ws.Range("L31").Value = 12528.68
$ws.Range("M31").Value = -3777.7222
$ws.Range("N31").Value = -13118.68
$ws.Range("H34").Value = 8988.977000000001
$ws.Range("I34").Value = 4072.7222
$ws.Range("J34").Value = 12528.68
$ws.Range("K34").Value = 4072.7222
$ws.Range("L34").Value = 12528.68
$ws.Range("M34").Value = -3870.7222
$ws.Range("N34").Value = -12932.68
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H86").Value = 6396456
$ws.Range("I86").Value = 7107154.5
$ws.Range("K86").Value = 7107154.5
$ws.Range("M86").Value = -7106031.5
$ws.Range("H89").Value = 6396456
$ws.Range("I89").Value = 7107154.5
$ws.Range("K89").Value = 35535772.5
$ws.Range("M89").Value = -35530156.5
$ws.Range("H100").Value = 39962.332
$ws.Range("J100").Value = 39962.332
$ws.Range("L100").Value = 39962.332
$ws.Range("N100").Value = -42126.332
$ws.Range("H141").Value = 56759.285
$ws.Range("J141").Value = 56759.285
$ws.Range("L141").Value = 56759.285
$ws.Range("N141").Value = -67119.285
$ws.Range("N50").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 9617009
$ws.Range("J92").Value = 9617009
$ws.Range("L92").Value = 28851027
$ws.Range("N92").Value = -28853523
$ws.Range("H107").Value = 22222520
$ws.Range("J107").Value = 28571700
$ws.Range("L107").Value = 85715100
$ws.Range("N107").Value = -85718940
$ws.Range("H113").Value = 5185.5
$ws.Range("J113").Value = 8119.6
$ws.Range("L113").Value = 24358.8
$ws.Range("N113").Value = -28698.8
$ws.Range("H128").Value = 219999.33
$ws.Range("I128").Value = 219999.33
$ws.Range("K128").Value = 659997.99
$ws.Range("M128").Value = -655017.99
$ws.Range("H134").Value = 80583.38
$ws.Range("I134").Value = 80583.38
$ws.Range("K134").Value = 241750.14
$ws.Range("M134").Value = -236680.14

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3638.1
$ws.Range("I80").Value = 2799.5
$ws.Range("J80").Value = 4197.1665
$ws.Range("K80").Value = 2799.5
$ws.Range("L80").Value = 4197.1665
$ws.Range("M80").Value = -1801.5
$ws.Range("N80").Value = -6193.1665
$ws.Range("H83").Value = 3638.1
$ws.Range("I83").Value = 2799.5
$ws.Range("J83").Value = 4197.1665
$ws.Range("K83").Value = 13997.5
$ws.Range("L83").Value = 20985.8325
$ws.Range("M83").Value = -9005.5
$ws.Range("N83").Value = -30969.8325
$ws.Range("H132").Value = 4808.375
$ws.Range("I132").Value = 1905.7222
$ws.Range("K132").Value = 5717.1666
$ws.Range("M132").Value = -3187.1666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 1599.5
$ws.Range("I107").Value = 1599.5
$ws.Range("K107").Value = 1599.5
$ws.Range("M107").Value = 320.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2234.5833
$ws.Range("I113").Value = 1464.909
$ws.Range("K113").Value = 4394.727000000001
$ws.Range("M113").Value = -2224.727000000001
$ws.Range("H132").Value = 22745652
$ws.Range("I132").Value = 45464624
$ws.Range("J132").Value = 26681.273
$ws.Range("K132").Value = 136393872
$ws.Range("L132").Value = 80043.819
$ws.Range("M132").Value = -136391342
$ws.Range("N132").Value = -85103.819
$ws.Range("H133").Value = 128743
$ws.Range("J133").Value = 128743
$ws.Range("L133").Value = 128743
$ws.Range("N133").Value = -138863
